# Generate Report for Handoff
# Rename the handed-off source file (and its generated handoff packages) from the
# old GUID to a new GUID, across the Overview / zh-cn / de-de sheets, updating both
# the cell text and the matching hyperlink display text. The hyperlink targets
# themselves are left exactly as they were.

$wb = $excel.ActiveWorkbook

$newGuidFile = "6cd60f7e-c9a3-481d-b9f8-7d00f8083620.md"
$newZhCnXlf  = "6cd60f7e-c9a3-481d-b9f8-7d00f8083620.009381db6f9cfc8d1d2dcac076c1b0a9cdac6e81.zh-cn.xlf"
$newDeDeXlf  = "6cd60f7e-c9a3-481d-b9f8-7d00f8083620.009381db6f9cfc8d1d2dcac076c1b0a9cdac6e81.de-de.xlf"

$newZhCnTime = "2016-02-22 17:54:33"
$newDeDeTime = "2016-02-22 17:54:44"

$configDisplay = ".localization-config"

$linkMdOverview     = "https://github.com/OpenLocalizationTest/oltest/blob/47859df1207a7278e194e60322f6f58b3cd2299e/e2e/7d4d153a-6cce-48b8-86ee-11b21f53d269.md"
$linkConfigOverview = "https://github.com/OpenLocalizationTest/oltest/blob/47859df1207a7278e194e60322f6f58b3cd2299e/.localization-config"

$linkMdZhCn     = "https://github.com/OpenLocalizationTest/oltest/blob/47859df1207a7278e194e60322f6f58b3cd2299e/e2e/7d4d153a-6cce-48b8-86ee-11b21f53d269.md"
$linkXlfZhCn    = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f9cc27dcd48cc94873a4ef7d7af47a049d2bd451/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/7d4d153a-6cce-48b8-86ee-11b21f53d269.eb26e5d07873e7698ed57e44c178012315358b19.zh-cn.xlf"
$linkConfigZhCn = "https://github.com/OpenLocalizationTest/oltest/blob/47859df1207a7278e194e60322f6f58b3cd2299e/.localization-config"

$linkMdDeDe     = "https://github.com/OpenLocalizationTest/oltest/blob/47859df1207a7278e194e60322f6f58b3cd2299e/e2e/7d4d153a-6cce-48b8-86ee-11b21f53d269.md"
$linkXlfDeDe    = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9cb8507a1ae1ac405f443efd27cd554ea9b58418/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/7d4d153a-6cce-48b8-86ee-11b21f53d269.eb26e5d07873e7698ed57e44c178012315358b19.de-de.xlf"
$linkConfigDeDe = "https://github.com/OpenLocalizationTest/oltest/blob/47859df1207a7278e194e60322f6f58b3cd2299e/.localization-config"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value2 = $newGuidFile

$wsOverview.Range("A2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $linkMdOverview, "", "", $newGuidFile) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $linkConfigOverview, "", "", $configDisplay) | Out-Null

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value2 = $newGuidFile
$wsZhCn.Range("C2").Value2 = $newZhCnXlf
$wsZhCn.Range("D2").Value2 = $newZhCnTime

$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $linkMdZhCn, "", "", $newGuidFile) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C2"), $linkXlfZhCn, "", "", $newZhCnXlf) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $linkConfigZhCn, "", "", $configDisplay) | Out-Null

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value2 = $newGuidFile
$wsDeDe.Range("C2").Value2 = $newDeDeXlf
$wsDeDe.Range("D2").Value2 = $newDeDeTime

$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $linkMdDeDe, "", "", $newGuidFile) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C2"), $linkXlfDeDe, "", "", $newDeDeXlf) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $linkConfigDeDe, "", "", $configDisplay) | Out-Null
